$wb = $excel.ActiveWorkbook

# --- 1. Status text update: "Ready for handoff" -> "In Translation" -------
# Overview sheet: columns E (zh-cn) and F (de-de), row 2
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# Per-language report sheets: column C ("Status"), row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2. Shrink the "Status" columns to fit the shorter text ---------------
# The shorter status string makes AutoFit pick a narrower column; COM's
# ColumnWidth setter snaps to the host's pixel grid, so we pick the input
# that lands on the grid point nearest the fitted width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
